$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 420
$ws.Range("I29").Value = 420
$ws.Range("K29").Value = 1260
$ws.Range("M29").Value = -979

$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516

$ws.Range("H112").Value = 1823.5312
$ws.Range("J112").Value = 1894.931
$ws.Range("L112").Value = 5684.793
$ws.Range("N112").Value = -7900.793

$ws.Range("H137").Value = 4018.2104
$ws.Range("I137").Value = 3552.7334
$ws.Range("K137").Value = 10658.2002
$ws.Range("M137").Value = -8108.200199999999

$ws.Range("H138").Value = 5053444.5
$ws.Range("I138").Value = 1158.3334
$ws.Range("J138").Value = 7411177.5
$ws.Range("K138").Value = 3475.0002
$ws.Range("L138").Value = 22233532.5
$ws.Range("M138").Value = 1664.9998
$ws.Range("N138").Value = -22243812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28860.818
$ws.Range("I2").Value = 50726.332
$ws.Range("K2").Value = 50726.332
$ws.Range("M2").Value = -50613.332

$ws.Range("H32").Value = 11373.143
$ws.Range("I32").Value = 9126.266
$ws.Range("J32").Value = 35339.832
$ws.Range("K32").Value = 9126.266
$ws.Range("L32").Value = 35339.832
$ws.Range("M32").Value = -8839.266
$ws.Range("N32").Value = -35913.832

$ws.Range("H61").Value = 4360.4
$ws.Range("I61").Value = 4815.7144
$ws.Range("K61").Value = 4815.7144
$ws.Range("M61").Value = -4603.7144

$ws.Range("H74").Value = 2252.276
$ws.Range("I74").Value = 2279.1428
$ws.Range("K74").Value = 2279.1428
$ws.Range("M74").Value = -1405.1428

$ws.Range("H77").Value = 2252.276
$ws.Range("I77").Value = 2279.1428
$ws.Range("K77").Value = 11395.714
$ws.Range("M77").Value = -7027.714

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H116").Value = 28860.818
$ws.Range("I116").Value = 50726.332
$ws.Range("K116").Value = 50726.332
$ws.Range("M116").Value = -48432.332

$ws.Range("H122").Value = 4149.077
$ws.Range("J122").Value = 10447.444
$ws.Range("L122").Value = 31342.332
$ws.Range("N122").Value = -36242.33199999999

$ws.Range("H132").Value = 4292.5684
$ws.Range("I132").Value = 3958.7297
$ws.Range("K132").Value = 11876.1891
$ws.Range("M132").Value = -9346.1891

$ws.Range("H136").Value = 4360.4
$ws.Range("I136").Value = 4815.7144
$ws.Range("K136").Value = 14447.1432
$ws.Range("M136").Value = -11897.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28860.818
$ws.Range("I3").Value = 50726.332
$ws.Range("K3").Value = 50726.332
$ws.Range("M3").Value = -50612.332

$ws.Range("H22").Value = 3003.4285
$ws.Range("I22").Value = 4135
$ws.Range("J22").Value = 174.5
$ws.Range("K22").Value = 4135
$ws.Range("L22").Value = 174.5
$ws.Range("M22").Value = -3962
$ws.Range("N22").Value = -520.5

$ws.Range("H86").Value = 3766.2693
$ws.Range("J86").Value = 5384.4614
$ws.Range("L86").Value = 5384.4614
$ws.Range("N86").Value = -7630.4614

$ws.Range("H89").Value = 3766.2693
$ws.Range("J89").Value = 5384.4614
$ws.Range("L89").Value = 26922.307
$ws.Range("N89").Value = -38154.307

$ws.Range("H94").Value = 6616.6665
$ws.Range("I94").Value = 5347.5
$ws.Range("J94").Value = 9155
$ws.Range("K94").Value = 5347.5
$ws.Range("L94").Value = 9155
$ws.Range("M94").Value = -4896.5
$ws.Range("N94").Value = -10057

$ws.Range("H96").Value = 49999.5
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 98999
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 98999
$ws.Range("M96").Value = 1746
$ws.Range("N96").Value = -104491

$ws.Range("H109").Value = 109999.6
$ws.Range("J109").Value = 109999.6
$ws.Range("L109").Value = 109999.6
$ws.Range("N109").Value = -112773.6

$ws.Range("H134").Value = 2813.3125
$ws.Range("I134").Value = 1867.3077
$ws.Range("J134").Value = 6912.6665
$ws.Range("K134").Value = 5601.9231
$ws.Range("L134").Value = 20737.9995
$ws.Range("M134").Value = -3066.9231
$ws.Range("N134").Value = -25807.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 625026.25
$ws.Range("I4").Value = 1000024
$ws.Range("K4").Value = 1000024
$ws.Range("M4").Value = -999912

$ws.Range("H31").Value = 2742.7942
$ws.Range("I31").Value = 2015.3334
$ws.Range("K31").Value = 2015.3334
$ws.Range("M31").Value = -1720.3334

$ws.Range("H34").Value = 2742.7942
$ws.Range("I34").Value = 2015.3334
$ws.Range("K34").Value = 2015.3334
$ws.Range("M34").Value = -1813.3334

$ws.Range("H122").Value = 331930
$ws.Range("I122").Value = 537976.5
$ws.Range("J122").Value = 5689.75
$ws.Range("K122").Value = 1613929.5
$ws.Range("L122").Value = 17069.25
$ws.Range("M122").Value = -1611479.5
$ws.Range("N122").Value = -21969.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4435.4443
$ws.Range("I3").Value = 2988.4285
$ws.Range("J3").Value = 9500
$ws.Range("K3").Value = 8965.2855
$ws.Range("L3").Value = 28500
$ws.Range("M3").Value = -8853.2855
$ws.Range("N3").Value = -28724

$ws.Range("H136").Value = 1861.5555
$ws.Range("J136").Value = 4075
$ws.Range("L136").Value = 12225
$ws.Range("N136").Value = -22425

$ws.Range("H139").Value = 45458096
$ws.Range("I139").Value = 50003056
$ws.Range("K139").Value = 150009168
$ws.Range("M139").Value = -150004028

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2742.1667
$ws.Range("I102").Value = 2014.2727
$ws.Range("K102").Value = 2014.2727
$ws.Range("M102").Value = -392.2727

$ws.Range("H113").Value = 7609.357
$ws.Range("I113").Value = 4282.1113
$ws.Range("K113").Value = 4282.1113
$ws.Range("M113").Value = -2112.1113

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 3537.8667
$ws.Range("I132").Value = 3004.88
$ws.Range("K132").Value = 9014.639999999999
$ws.Range("M132").Value = -6484.639999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2248.111
$ws.Range("I61").Value = 1152.9642
$ws.Range("K61").Value = 1152.9642
$ws.Range("M61").Value = -950.9641999999999

$ws.Range("H113").Value = 2248.111
$ws.Range("I113").Value = 1152.9642
$ws.Range("K113").Value = 1152.9642
$ws.Range("M113").Value = 1017.0358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H136").Value = 2746.5715
$ws.Range("I136").Value = 1874.4572
$ws.Range("J136").Value = 7107.143
$ws.Range("K136").Value = 5623.3716
$ws.Range("L136").Value = 21321.429
$ws.Range("M136").Value = -3073.3716
$ws.Range("N136").Value = -26421.429
